# "Actualizacion desde MV -datos-": append the new "04-10-2021" reading as
# row 68 of Sheet1 (same columns/shape as the existing daily rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds date-like labels stored as plain text (shared strings),
# not real dates. Assigning the string "04-10-2021" straight to Range.Value
# would be auto-recognized by Excel as a date and stored as a date serial
# number with a date number format. To get a genuine text value instead,
# write it as a text formula in a scratch cell, copy it, and paste only the
# resulting value into A68 - that preserves it as text.
$ws.Range("Z1").Formula = "=""04-10-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A68").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").ClearContents()

$ws.Range("B68").Value = 50000
$ws.Range("C68").Value = 45000
$ws.Range("D68").Value = 45000
$ws.Range("E68").Value = 40000
$ws.Range("F68").Value = 5000
$ws.Range("G68").Value = 2.25
